# Ajustes para que se guarde en el escritorio
# Adds two new survey entries (rows 7 and 8) to the registros sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Roger Villegas
$ws.Range("A7").Value = "Roger "
$ws.Range("B7").Value = "Villegas "
$ws.Range("C7").Value = "natación"
$ws.Range("D7").Value = "Masculino"
$ws.Range("E7").Value = "Santa Rosa"
$ws.Range("F7").Value = "Sí"
$ws.Range("G7").Value = "ford, nissan"

# Row 8: Ruben Tuesta
$ws.Range("A8").Value = "Ruben "
$ws.Range("B8").Value = "Tuesta "
$ws.Range("C8").Value = "béisbol"
$ws.Range("D8").Value = "Masculino"
$ws.Range("E8").Value = "Huehuetenango"
$ws.Range("F8").Value = "No"
$ws.Range("G8").Value = "ford, nissan"
